$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5368.4546
$ws.Range("I64").Value = 4065.6667
$ws.Range("K64").Value = 4065.6667
$ws.Range("M64").Value = -3817.6667

$ws.Range("H67").Value = 5368.4546
$ws.Range("I67").Value = 4065.6667
$ws.Range("K67").Value = 4065.6667
$ws.Range("M67").Value = -3207.6667

$ws.Range("H70").Value = 1525.5555
$ws.Range("J70").Value = 1409.75
$ws.Range("L70").Value = 4229.25
$ws.Range("N70").Value = -4769.25

$ws.Range("H73").Value = 1525.5555
$ws.Range("J73").Value = 1409.75
$ws.Range("L73").Value = 4229.25
$ws.Range("N73").Value = -6101.25

$ws.Range("H107").Value = 6829.1562
$ws.Range("I107").Value = 3349.3914
$ws.Range("K107").Value = 3349.3914
$ws.Range("M107").Value = -1429.3914

$ws.Range("H137").Value = 16525.525
$ws.Range("I137").Value = 8666.5
$ws.Range("K137").Value = 25999.5
$ws.Range("M137").Value = -23449.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4061.1914
$ws.Range("I32").Value = 4061.1914
$ws.Range("K32").Value = 4061.1914
$ws.Range("M32").Value = -3774.1914

$ws.Range("H122").Value = 1663.0526
$ws.Range("I122").Value = 1520.6562
$ws.Range("K122").Value = 4561.9686
$ws.Range("M122").Value = -2111.9686

$ws.Range("H132").Value = 3069.6
$ws.Range("I132").Value = 2931.7856
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8795.356800000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6265.356800000001
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 9998
$ws.Range("I29").Value = 9998
$ws.Range("K29").Value = 9998
$ws.Range("M29").Value = -9709

$ws.Range("H32").Value = 25929
$ws.Range("J32").Value = 25929
$ws.Range("L32").Value = 25929
$ws.Range("N32").Value = -26697

$ws.Range("H34").Value = 5750
$ws.Range("J34").Value = 5750
$ws.Range("L34").Value = 5750
$ws.Range("N34").Value = -5978

$ws.Range("H107").Value = 903.78845
$ws.Range("I107").Value = 642.83673
$ws.Range("J107").Value = 5166
$ws.Range("K107").Value = 642.83673
$ws.Range("L107").Value = 5166
$ws.Range("M107").Value = 1277.16327
$ws.Range("N107").Value = -9006

$ws.Range("H134").Value = 8681.298000000001
$ws.Range("I134").Value = 4333.9
$ws.Range("J134").Value = 16353.177
$ws.Range("K134").Value = 13001.7
$ws.Range("L134").Value = 49059.531
$ws.Range("M134").Value = -10466.7
$ws.Range("N134").Value = -54129.531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 39116.125
$ws.Range("J23").Value = 39134.285
$ws.Range("L23").Value = 39134.285
$ws.Range("N23").Value = -39614.285

$ws.Range("H27").Value = 39116.125
$ws.Range("J27").Value = 39134.285
$ws.Range("L27").Value = 39134.285
$ws.Range("N27").Value = -39518.285

$ws.Range("H31").Value = 2539.0833
$ws.Range("J31").Value = 3489.2307
$ws.Range("L31").Value = 3489.2307
$ws.Range("N31").Value = -4079.2307

$ws.Range("H34").Value = 2539.0833
$ws.Range("J34").Value = 3489.2307
$ws.Range("L34").Value = 3489.2307
$ws.Range("N34").Value = -3893.2307

$ws.Range("H58").Value = 11011.223
$ws.Range("I58").Value = 7495
$ws.Range("K58").Value = 7495
$ws.Range("M58").Value = -7292

$ws.Range("H134").Value = 6213.643
$ws.Range("I134").Value = 5064.4707
$ws.Range("K134").Value = 15193.4121
$ws.Range("M134").Value = -12658.4121

$ws.Range("H136").Value = 11011.223
$ws.Range("I136").Value = 7495
$ws.Range("K136").Value = 22485
$ws.Range("M136").Value = -19935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2159.5715
$ws.Range("I5").Value = 1736
$ws.Range("J5").Value = 2275.0908
$ws.Range("K5").Value = 5208
$ws.Range("L5").Value = 6825.2724
$ws.Range("M5").Value = -5096
$ws.Range("N5").Value = -7049.2724

$ws.Range("H135").Value = 2159.5715
$ws.Range("I135").Value = 1736
$ws.Range("J135").Value = 2275.0908
$ws.Range("K135").Value = 15624
$ws.Range("L135").Value = 20475.8172
$ws.Range("M135").Value = -13089
$ws.Range("N135").Value = -25545.8172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1089.1666
$ws.Range("I97").Value = 973
$ws.Range("K97").Value = 973
$ws.Range("M97").Value = -477

$ws.Range("H102").Value = 1267.742
$ws.Range("I102").Value = 1251.7693
$ws.Range("K102").Value = 1251.7693
$ws.Range("M102").Value = 370.2307000000001

$ws.Range("H113").Value = 235340.56
$ws.Range("I113").Value = 252258.12
$ws.Range("J113").Value = 100000
$ws.Range("K113").Value = 252258.12
$ws.Range("L113").Value = 100000
$ws.Range("M113").Value = -250088.12
$ws.Range("N113").Value = -104340

$ws.Range("H122").Value = 1750.5652
$ws.Range("I122").Value = 1705.8422
$ws.Range("J122").Value = 1963
$ws.Range("K122").Value = 5117.5266
$ws.Range("L122").Value = 5889
$ws.Range("M122").Value = -2667.5266
$ws.Range("N122").Value = -10789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9039.5
$ws.Range("I7").Value = 7513.091
$ws.Range("K7").Value = 7513.091
$ws.Range("M7").Value = -7401.091

$ws.Range("H93").Value = 2735
$ws.Range("J93").Value = 1433.8572
$ws.Range("L93").Value = 1433.8572
$ws.Range("N93").Value = -3929.8572

$ws.Range("H126").Value = 9039.5
$ws.Range("I126").Value = 7513.091
$ws.Range("K126").Value = 22539.273
$ws.Range("M126").Value = -20069.273

$ws.Range("H136").Value = 7711.55
$ws.Range("I136").Value = 7469.1
$ws.Range("J136").Value = 7954
$ws.Range("K136").Value = 22407.3
$ws.Range("L136").Value = 23862
$ws.Range("M136").Value = -19857.3
$ws.Range("N136").Value = -28962

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H132").Value = 222094.7
$ws.Range("I132").Value = 336361.56
$ws.Range("J132").Value = 35659.316
$ws.Range("K132").Value = 1009084.68
$ws.Range("L132").Value = 106977.948
$ws.Range("M132").Value = -1006554.68
$ws.Range("N132").Value = -112037.948
